$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.43"
$ws.Range("E2").Value = "'-0.79%"
$ws.Range("D3").Value = "'35.45"
$ws.Range("E3").Value = "'3.86%"
$ws.Range("D4").Value = "'5.045"
$ws.Range("E4").Value = "'-2.26%"
$ws.Range("D5").Value = "'0.07691"
$ws.Range("E5").Value = "'-1.83%"
$ws.Range("D6").Value = "'2.176"
$ws.Range("E6").Value = "'-9.79%"
$ws.Range("D7").Value = "'8.027"
$ws.Range("E7").Value = "'-0.34%"
$ws.Range("D8").Value = "'4.007"
$ws.Range("E8").Value = "'2.63%"
$ws.Range("D9").Value = "'0.9296"
$ws.Range("E9").Value = "'-0.65%"
$ws.Range("D10").Value = "'0.09343"
$ws.Range("E10").Value = "'-4.86%"
$ws.Range("D11").Value = "'0.1825"
$ws.Range("E11").Value = "'2.56%"
$ws.Range("D12").Value = "'0.08483"
$ws.Range("E12").Value = "'-0.61%"
$ws.Range("D13").Value = "'0.03598"
$ws.Range("E13").Value = "'7.26%"
$ws.Range("D14").Value = "'0.09982"
$ws.Range("E14").Value = "'0.59%"
$ws.Range("D15").Value = "'0.001485"
$ws.Range("E15").Value = "'0.30%"
$ws.Range("D16").Value = "'0.005757"
$ws.Range("E16").Value = "'1.22%"
$ws.Range("D17").Value = "'3.478"
$ws.Range("E17").Value = "'0.35%"
$ws.Range("D18").Value = "'2.184"
$ws.Range("E18").Value = "'1.01%"
$ws.Range("D19").Value = "'0.3463"
$ws.Range("E19").Value = "'2.86%"
$ws.Range("D20").Value = "'0.1328"
$ws.Range("E20").Value = "'-1.50%"
$ws.Range("D21").Value = "'4.582"
$ws.Range("E21").Value = "'6.86%"
$ws.Range("E22").Value = "'-2.02%"
$ws.Range("D23").Value = "'0.04665"
$ws.Range("E23").Value = "'0.36%"
$ws.Range("D24").Value = "'0.001239"
$ws.Range("E24").Value = "'1.50%"
$ws.Range("D25").Value = "'0.004474"
$ws.Range("E25").Value = "'1.28%"
$ws.Range("E26").Value = "'0.79%"
$ws.Range("E27").Value = "'-20.23%"
$ws.Range("D39").Value = "'0.01723"
$ws.Range("E39").Value = "'-1.27%"
$ws.Range("D40").Value = "'0.04680"
$ws.Range("E40").Value = "'-2.80%"
$ws.Range("D41").Value = "'0.007947"
$ws.Range("E41").Value = "'2.44%"
$ws.Range("D42").Value = "'0.1400"
$ws.Range("E42").Value = "'-1.06%"
$ws.Range("D43").Value = "'0.007696"
$ws.Range("E43").Value = "'-21.50%"
$ws.Range("D44").Value = "'0.002234"
$ws.Range("E44").Value = "'7.41%"
$ws.Range("D45").Value = "'0.008949"
$ws.Range("E45").Value = "'-2.03%"
$ws.Range("D46").Value = "'0.00006251"
$ws.Range("E46").Value = "'2.63%"
$ws.Range("D47").Value = "'0.00000000756"
$ws.Range("E47").Value = "'1.03%"
$ws.Range("D48").Value = "'5.736"
$ws.Range("E48").Value = "'116.08%"
$ws.Range("D49").Value = "'0.002711"
$ws.Range("E49").Value = "'35.91%"
$ws.Range("D50").Value = "'0.00002116"
$ws.Range("E50").Value = "'1.03%"
$ws.Range("D51").Value = "'0.0002015"
$ws.Range("E51").Value = "'1.03%"
